# Tried a bubble sort.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31 (Big-O) - Practice 1 (J31) done: change from "Sat" (date/day) to numeric 4
$ws.Range("J31").Value = 4

# Row 33 (Ordered List)
$ws.Range("H33").Value = "Must keep the order; remember to include all methods (index, search, size, etc.)"
$ws.Range("J33").Value = 4

# Row 34 (Bubble Sort)
$ws.Range("F34").Value = 1
$ws.Range("G34").Value = "n^2"
$ws.Range("H34").Value = "Extremely slow, worst of the sorts, but fine for almost-sorted lists; easiest implementation for small lists too"
$ws.Range("J34").Value = 4

# Update the active selection to reflect where the user ended up (J35)
$ws.Range("J35").Select()
